# Apply the diff changes to both the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row for "张家港· 喵喵漫国潮动漫节" (F5: 46 -> 48)
    $ws.Range("F5").Value = 48

    # Row for "苏州·漫遇引力动漫游戏展" (F6: 607 -> 608)
    $ws.Range("F6").Value = 608

    # Row for "苏州·首届 童年回忆同人only 茶歇聚会" (row 7)
    # Name gains a "（取消）" suffix, and the lowest-price column switches
    # from a numeric value to the text "不可售".
    $ws.Range("C7").Value = "苏州·首届 童年回忆同人only 茶歇聚会（取消）"
    $ws.Range("G7").Value = "不可售"
}

# "展览" sheet specific row numbers for the remaining F-column updates.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F12").Value = 102
$wsExhibit.Range("F13").Value = 303
$wsExhibit.Range("F17").Value = 11171
$wsExhibit.Range("F18").Value = 5336

# "全部类型" sheet specific row numbers for the remaining F-column updates.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 102
$wsAll.Range("F15").Value = 303
$wsAll.Range("F19").Value = 11171
$wsAll.Range("F21").Value = 5336
